$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Save the original row 5 content before overwriting it
$origB5 = $ws.Cells.Item(5, 2).Value()
$origC5 = $ws.Cells.Item(5, 3).Value()
$origD5 = $ws.Cells.Item(5, 4).Value()

# Update D3/D4 from "no" to "yes"
$ws.Cells.Item(3, 4).Value = "yes"
$ws.Cells.Item(4, 4).Value = "yes"

# Update row 5 (TC_EC_0004) with the new negative-scenario description
$ws.Cells.Item(5, 2).Value = "fetch the particular Pet in Petstore swagger which is not in database"
$ws.Cells.Item(5, 3).Value = "Regression"
$ws.Cells.Item(5, 4).Value = "yes"

# Add a new row 6 that carries what used to be in row 5 (new Test Case ID)
$ws.Cells.Item(6, 1).Value = "TC_EC_0005"
$ws.Cells.Item(6, 2).Value = $origB5
$ws.Cells.Item(6, 3).Value = $origC5
$ws.Cells.Item(6, 4).Value = $origD5

$ws.Range("D6").Select()
